$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header: "Expected" (C1)
$ws.Range("C1").Value = "Expected"

# New column C, row 2: a hyperlink to the book's Goodreads page.
# Excel automatically uses the URL as the cell's displayed text and
# applies the built-in "Hyperlink" style (underlined, themed font) -
# this also creates the "Hyperlink" font/cellStyle entries in styles.xml
# and the external relationship the <hyperlinks> element points at.
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.goodreads.com/book/show/865.The_Alchemist")

# Size the columns to fit their new contents (author column A and the
# new link column C).
$ws.Columns("A").ColumnWidth = 11.833333333333334
$ws.Columns("C").ColumnWidth = 55.833333333333336

# Leave the selection where data entry would continue next.
$ws.Range("C3").Select() | Out-Null
